# Auto-generated edit script applying numeric corrections to the
# "Seraph_Profits" workbook (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block updates the H:N price/profit columns for a specific leve row,
# matching the authoritative values from the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row6
$ws.Range("H6").Value = 189.08333
$ws.Range("I6").Value = 196.9
$ws.Range("K6").Value = 590.7
$ws.Range("M6").Value = -478.7

# ALC!row18
$ws.Range("H18").Value = 1374.5
$ws.Range("I18").Value = 1374.5
$ws.Range("K18").Value = 1374.5
$ws.Range("M18").Value = -1090.5

# ALC!row43
$ws.Range("H43").Value = 6135.2354
$ws.Range("I43").Value = 4849.75
$ws.Range("K43").Value = 4849.75
$ws.Range("M43").Value = -4780.75

# ALC!row58
$ws.Range("H58").Value = 2049
$ws.Range("J58").Value = 2537.5
$ws.Range("L58").Value = 7612.5
$ws.Range("N58").Value = -7912.5

# ALC!row76
$ws.Range("H76").Value = 2805.25
$ws.Range("I76").Value = 2740.3333
$ws.Range("K76").Value = 2740.3333
$ws.Range("M76").Value = -2425.3333

# ALC!row79
$ws.Range("H79").Value = 2805.25
$ws.Range("I79").Value = 2740.3333
$ws.Range("K79").Value = 2740.3333
$ws.Range("M79").Value = -1648.3333

# ALC!row92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# ALC!row116
$ws.Range("H116").Value = 4921.75
$ws.Range("I116").Value = 4900
$ws.Range("K116").Value = 4900
$ws.Range("M116").Value = -1458

# ALC!row132
$ws.Range("H132").Value = 1436.7778
$ws.Range("I132").Value = 1541.5
$ws.Range("K132").Value = 4624.5
$ws.Range("M132").Value = -2094.5

# ALC!row138
$ws.Range("H138").Value = 2651.8333
$ws.Range("J138").Value = 3569.0476
$ws.Range("L138").Value = 10707.1428
$ws.Range("N138").Value = -20987.1428

$ws = $wb.Worksheets.Item("ARM")
# ARM!row12
$ws.Range("H12").Value = 600979.8
$ws.Range("J12").Value = 1299.6666
$ws.Range("L12").Value = 1299.6666
$ws.Range("N12").Value = -1645.6666

# ARM!row37
$ws.Range("H37").Value = 15408.818
$ws.Range("J37").Value = 21356.715
$ws.Range("L37").Value = 21356.715
$ws.Range("N37").Value = -21902.715

# ARM!row44
$ws.Range("H44").Value = 14993.667
$ws.Range("J44").Value = 14993.667
$ws.Range("L44").Value = 14993.667
$ws.Range("N44").Value = -15969.667

# ARM!row55
$ws.Range("H55").Value = 25179.4
$ws.Range("J55").Value = 33949.668
$ws.Range("L55").Value = 33949.668
$ws.Range("N55").Value = -34579.668

$ws = $wb.Worksheets.Item("BSM")
# BSM!row38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

# BSM!row99
$ws.Range("H99").Value = 1711.3334
$ws.Range("I99").Value = 1408.1177
$ws.Range("K99").Value = 1408.1177
$ws.Range("M99").Value = 89.88229999999999

# BSM!row100
$ws.Range("H100").Value = 19528
$ws.Range("J100").Value = 19528
$ws.Range("L100").Value = 19528
$ws.Range("N100").Value = -21692

$ws = $wb.Worksheets.Item("CRP")
# CRP!row7
$ws.Range("H7").Value = 146
$ws.Range("I7").Value = 104.72
$ws.Range("J7").Value = 352.4
$ws.Range("K7").Value = 104.72
$ws.Range("L7").Value = 352.4
$ws.Range("M7").Value = 8.280000000000001
$ws.Range("N7").Value = -578.4

# CRP!row15
$ws.Range("H15").Value = 15615.667
$ws.Range("I15").Value = 14419.5
$ws.Range("J15").Value = 18008
$ws.Range("K15").Value = 14419.5
$ws.Range("L15").Value = 18008
$ws.Range("M15").Value = -14249.5
$ws.Range("N15").Value = -18348

# CRP!row22
$ws.Range("H22").Value = 736.8570999999999
$ws.Range("I22").Value = 736.8570999999999
$ws.Range("K22").Value = 736.8570999999999
$ws.Range("M22").Value = -386.8570999999999

# CRP!row29
$ws.Range("H29").Value = 6077
$ws.Range("I29").Value = 6000
$ws.Range("J29").Value = 6115.5
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 6115.5
$ws.Range("M29").Value = -5707
$ws.Range("N29").Value = -6701.5

# CRP!row31
$ws.Range("H31").Value = 3659.9167
$ws.Range("I31").Value = 2436.5557
$ws.Range("K31").Value = 2436.5557
$ws.Range("M31").Value = -2141.5557

# CRP!row34
$ws.Range("H34").Value = 3659.9167
$ws.Range("I34").Value = 2436.5557
$ws.Range("K34").Value = 2436.5557
$ws.Range("M34").Value = -2234.5557

# CRP!row54
$ws.Range("H54").Value = 20091
$ws.Range("J54").Value = 20091
$ws.Range("L54").Value = 20091
$ws.Range("N54").Value = -21407

# CRP!row86
$ws.Range("H86").Value = 8249
$ws.Range("I86").Value = 7000
$ws.Range("J86").Value = 8665.333000000001
$ws.Range("K86").Value = 7000
$ws.Range("L86").Value = 8665.333000000001
$ws.Range("M86").Value = -5877
$ws.Range("N86").Value = -10911.333

# CRP!row89
$ws.Range("H89").Value = 8249
$ws.Range("I89").Value = 7000
$ws.Range("J89").Value = 8665.333000000001
$ws.Range("K89").Value = 35000
$ws.Range("L89").Value = 43326.665
$ws.Range("M89").Value = -29384
$ws.Range("N89").Value = -54558.665

$ws = $wb.Worksheets.Item("CUL")
# CUL!row123
$ws.Range("H123").Value = 1000
$ws.Range("J123").Value = 1000
$ws.Range("L123").Value = 3000
$ws.Range("N123").Value = -7900

# CUL!row126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# CUL!row130
$ws.Range("H130").Value = 4494.3335
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()

# CUL!row131
$ws.Range("H131").Value = 995
$ws.Range("J131").Value = 995
$ws.Range("L131").Value = 2985
$ws.Range("N131").Value = -13065

$ws = $wb.Worksheets.Item("GSM")
# GSM!row6
$ws.Range("H6").Value = 400
$ws.Range("I6").Value = 400
$ws.Range("K6").Value = 400
$ws.Range("M6").Value = -287

# GSM!row16
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 400
$ws.Range("K16").Value = 400
$ws.Range("M16").Value = -150

# GSM!row92
$ws.Range("H92").Value = 14083.667
$ws.Range("J92").Value = 16125.5
$ws.Range("L92").Value = 16125.5
$ws.Range("N92").Value = -19869.5

# GSM!row126
$ws.Range("H126").Value = 4994.6665
$ws.Range("I126").Value = 4994.6665
$ws.Range("K126").Value = 14983.9995
$ws.Range("M126").Value = -12513.9995

$ws = $wb.Worksheets.Item("LTW")
# LTW!row7
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

# LTW!row19
$ws.Range("H19").Value = 6167.6665
$ws.Range("I19").Value = 1751.5
$ws.Range("K19").Value = 1751.5
$ws.Range("M19").Value = -1581.5

# LTW!row22
$ws.Range("H22").Value = 2242.6667
$ws.Range("I22").Value = 1511.125
$ws.Range("J22").Value = 3078.7144
$ws.Range("K22").Value = 1511.125
$ws.Range("L22").Value = 3078.7144
$ws.Range("M22").Value = -1216.125
$ws.Range("N22").Value = -3668.7144

# LTW!row27
$ws.Range("H27").Value = 2242.6667
$ws.Range("I27").Value = 1511.125
$ws.Range("J27").Value = 3078.7144
$ws.Range("K27").Value = 1511.125
$ws.Range("L27").Value = 3078.7144
$ws.Range("M27").Value = -1404.125
$ws.Range("N27").Value = -3292.7144

# LTW!row34
$ws.Range("H34").Value = 13498.75
$ws.Range("J34").Value = 7999.5
$ws.Range("L34").Value = 7999.5
$ws.Range("N34").Value = -8343.5

# LTW!row55
$ws.Range("H55").Value = 577.4
$ws.Range("I55").Value = 311
$ws.Range("J55").Value = 1199
$ws.Range("K55").Value = 311
$ws.Range("L55").Value = 1199
$ws.Range("M55").Value = -138
$ws.Range("N55").Value = -1545

# LTW!row126
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530

$ws = $wb.Worksheets.Item("WVR")
# WVR!row11
$ws.Range("H11").Value = 19992
$ws.Range("J11").Value = 19992
$ws.Range("L11").Value = 19992
$ws.Range("N11").Value = -20276

# WVR!row14
$ws.Range("H14").Value = 3619.7273
$ws.Range("I14").Value = 2363.2856
$ws.Range("K14").Value = 2363.2856
$ws.Range("M14").Value = -2195.2856

# WVR!row96
$ws.Range("H96").Value = 6501.3335
$ws.Range("J96").Value = 6501.3335
$ws.Range("L96").Value = 6501.3335
$ws.Range("N96").Value = -9247.333500000001

# WVR!row126
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530

# WVR!row132
$ws.Range("H132").Value = 2662.2104
$ws.Range("I132").Value = 2118.9333
$ws.Range("J132").Value = 4699.5
$ws.Range("K132").Value = 6356.7999
$ws.Range("L132").Value = 14098.5
$ws.Range("M132").Value = -3826.7999
$ws.Range("N132").Value = -19158.5

# WVR!row136
$ws.Range("H136").Value = 1890
$ws.Range("I136").Value = 1186.6666
$ws.Range("K136").Value = 3559.9998
$ws.Range("M136").Value = -1009.9998

